$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.114.45"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.26%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.964.23"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.95%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "379.91"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.33"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.57%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.545"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.95%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.592"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.74%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.57"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.40%  "

$ws.Range("E11").Value = "  -1.23%  "

$ws.Range("E12").Value = "  +2.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.427.65"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.82"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +6.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "18.35"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "11.99"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +67.48%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.960.99"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.69%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.00"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "51.189.20"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.12"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.45"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.96%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.57%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.18"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.79%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.80"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.47%  "

$ws.Range("E25").Value = "  +13.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.91"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.88%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.20"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -7.27%  "

$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.167"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.91"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.23%  "

$ws.Range("E31").Value = "  -1.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.41"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.77%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.49"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.92%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.07"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.53%  "

$ws.Range("E35").Value = "  +2.22%  "

$ws.Range("E36").Value = "  -2.63%  "

$ws.Range("E38").Value = "  +9.00%  "

$ws.Range("E39").Value = "  +1.87%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.84"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.88%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "16.59"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "125.21"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.11%  "

$ws.Range("E43").Value = "  -2.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.60"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.55"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +9.55%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.39"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.67%  "

$ws.Range("E47").Value = "  -0.54%  "

$ws.Range("B48").Value = "TheGraph"
$ws.Range("C48").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.272"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -5.93%  "

$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.053.58"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0321"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.89%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.42"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +7.74%  "
